# Update 苏州-漫展信息.xlsx output (gh-pages regeneration at 456a3b4)
# Applies updated "想去人数" (F) and "最低票价" (G) values to the
# "展览" sheet and the mirrored "全部类型" sheet.

$wb = $excel.ActiveWorkbook

function Update-Sheet {
    param($ws)

    # --- Column F ("想去人数") numeric updates ---
    $ws.Range("F2").Value  = 1223
    $ws.Range("F3").Value  = 14742
    $ws.Range("F4").Value  = 18070
    $ws.Range("F5").Value  = 18070
    $ws.Range("F7").Value  = 88
    $ws.Range("F16").Value = 66
    $ws.Range("F18").Value = 47
    $ws.Range("F19").Value = 1361

    # --- Column G ("最低票价") updates ---
    $ws.Range("G2").Value = 79.90000000000001

    $ws.Range("G8").Value  = "不可售"
    $ws.Range("G9").Value  = "不可售"
    $ws.Range("G10").Value = "不可售"
    $ws.Range("G11").Value = "不可售"
    $ws.Range("G20").Value = "不可售"
    $ws.Range("G21").Value = "不可售"
}

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
Update-Sheet $wsExhibit
# Row 17 differs only on this sheet (other sheet already had 178)
$wsExhibit.Range("F17").Value = 178
$wsExhibit.Range("F24").Value = 7470
$wsExhibit.Range("F28").Value = 1198
$wsExhibit.Range("F30").Value = 5898
$wsExhibit.Range("F31").Value = 78
$wsExhibit.Range("F36").Value = 5186

# --- Sheet "全部类型" (duplicated data, two extra rows vs "展览") ---
$wsAll = $wb.Worksheets.Item("全部类型")
Update-Sheet $wsAll
$wsAll.Range("F25").Value = 7470
$wsAll.Range("F29").Value = 1198
$wsAll.Range("F32").Value = 5898
$wsAll.Range("F33").Value = 78
$wsAll.Range("F38").Value = 5186
